# Auto-generated Excel COM-interop edit script
# Applies the BAJAJ-PL MIS Base Page data refresh to Sheet1 (rows 2-9, columns A-AE).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BKT=BKT0, STATE=DELHI NCR
$ws.Cells.Item(2, 1).Value = "BKT0"
$ws.Cells.Item(2, 2).Value = "DELHI NCR"
$ws.Cells.Item(2, 3).Value = 817120112.43
$ws.Cells.Item(2, 4).Value = 228
$ws.Cells.Item(2, 5).Value = 25
$ws.Cells.Item(2, 6).Value = 202
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 61842820.34
$ws.Cells.Item(2, 14).Value = 752412783.87
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 2864508.22
$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 7.57
$ws.Cells.Item(2, 22).Value = 92.08
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = 0
$ws.Cells.Item(2, 25).Value = 0.35
$ws.Cells.Item(2, 26).Value = 0
$ws.Cells.Item(2, 27).Value = 0
$ws.Cells.Item(2, 28).Value = 0
$ws.Cells.Item(2, 29).Value = 8915728
$ws.Cells.Item(2, 30).Value = 92.43
$ws.Cells.Item(2, 31).Value = 0.35

# Row 3: BKT=BKT1, STATE=DELHI NCR
$ws.Cells.Item(3, 1).Value = "BKT1"
$ws.Cells.Item(3, 2).Value = "DELHI NCR"
$ws.Cells.Item(3, 3).Value = 288503013.31
$ws.Cells.Item(3, 4).Value = 199
$ws.Cells.Item(3, 5).Value = 27
$ws.Cells.Item(3, 6).Value = 125
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 44
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 39913237.29
$ws.Cells.Item(3, 14).Value = 194592171.52
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 50315521.08
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 3682083.42
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 13.83
$ws.Cells.Item(3, 22).Value = 67.45
$ws.Cells.Item(3, 23).Value = 0
$ws.Cells.Item(3, 24).Value = 0
$ws.Cells.Item(3, 25).Value = 1.28
$ws.Cells.Item(3, 26).Value = 17.44
$ws.Cells.Item(3, 27).Value = 0
$ws.Cells.Item(3, 28).Value = 0
$ws.Cells.Item(3, 29).Value = 5234507
$ws.Cells.Item(3, 30).Value = 86.17
$ws.Cells.Item(3, 31).Value = 18.72

# Row 4: BKT=BKT1, STATE=RAJASTHAN
$ws.Cells.Item(4, 1).Value = "BKT1"
$ws.Cells.Item(4, 2).Value = "RAJASTHAN"
$ws.Cells.Item(4, 3).Value = 9911622.05
$ws.Cells.Item(4, 4).Value = 15
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 9911622.05
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 100
$ws.Cells.Item(4, 22).Value = 0
$ws.Cells.Item(4, 23).Value = 0
$ws.Cells.Item(4, 24).Value = 0
$ws.Cells.Item(4, 25).Value = 0
$ws.Cells.Item(4, 26).Value = 0
$ws.Cells.Item(4, 27).Value = 0
$ws.Cells.Item(4, 28).Value = 0
$ws.Cells.Item(4, 29).Value = 0
$ws.Cells.Item(4, 30).Value = 0
$ws.Cells.Item(4, 31).Value = 0

# Row 5: BKT=BKT2, STATE=RAJASTHAN
$ws.Cells.Item(5, 1).Value = "BKT2"
$ws.Cells.Item(5, 2).Value = "RAJASTHAN"
$ws.Cells.Item(5, 3).Value = 1921240
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 1921240
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 0
$ws.Cells.Item(5, 18).Value = 0
$ws.Cells.Item(5, 19).Value = 0
$ws.Cells.Item(5, 20).Value = 0
$ws.Cells.Item(5, 21).Value = 100
$ws.Cells.Item(5, 22).Value = 0
$ws.Cells.Item(5, 23).Value = 0
$ws.Cells.Item(5, 24).Value = 0
$ws.Cells.Item(5, 25).Value = 0
$ws.Cells.Item(5, 26).Value = 0
$ws.Cells.Item(5, 27).Value = 0
$ws.Cells.Item(5, 28).Value = 0
$ws.Cells.Item(5, 29).Value = 0
$ws.Cells.Item(5, 30).Value = 0
$ws.Cells.Item(5, 31).Value = 0

# Row 6: BKT=BKT3, STATE=RAJASTHAN
$ws.Cells.Item(6, 1).Value = "BKT3"
$ws.Cells.Item(6, 2).Value = "RAJASTHAN"
$ws.Cells.Item(6, 3).Value = 1765526.64
$ws.Cells.Item(6, 4).Value = 4
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = 1765526.64
$ws.Cells.Item(6, 14).Value = 0
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 0
$ws.Cells.Item(6, 18).Value = 0
$ws.Cells.Item(6, 19).Value = 0
$ws.Cells.Item(6, 20).Value = 0
$ws.Cells.Item(6, 21).Value = 100
$ws.Cells.Item(6, 22).Value = 0
$ws.Cells.Item(6, 23).Value = 0
$ws.Cells.Item(6, 24).Value = 0
$ws.Cells.Item(6, 25).Value = 0
$ws.Cells.Item(6, 26).Value = 0
$ws.Cells.Item(6, 27).Value = 0
$ws.Cells.Item(6, 28).Value = 0
$ws.Cells.Item(6, 29).Value = 0
$ws.Cells.Item(6, 30).Value = 0
$ws.Cells.Item(6, 31).Value = 0

# Row 7: BKT=BKT4, STATE=RAJASTHAN
$ws.Cells.Item(7, 1).Value = "BKT4"
$ws.Cells.Item(7, 2).Value = "RAJASTHAN"
$ws.Cells.Item(7, 3).Value = 489693
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 489693
$ws.Cells.Item(7, 14).Value = 0
$ws.Cells.Item(7, 15).Value = 0
$ws.Cells.Item(7, 16).Value = 0
$ws.Cells.Item(7, 17).Value = 0
$ws.Cells.Item(7, 18).Value = 0
$ws.Cells.Item(7, 19).Value = 0
$ws.Cells.Item(7, 20).Value = 0
$ws.Cells.Item(7, 21).Value = 100
$ws.Cells.Item(7, 22).Value = 0
$ws.Cells.Item(7, 23).Value = 0
$ws.Cells.Item(7, 24).Value = 0
$ws.Cells.Item(7, 25).Value = 0
$ws.Cells.Item(7, 26).Value = 0
$ws.Cells.Item(7, 27).Value = 0
$ws.Cells.Item(7, 28).Value = 0
$ws.Cells.Item(7, 29).Value = 0
$ws.Cells.Item(7, 30).Value = 0
$ws.Cells.Item(7, 31).Value = 0

# Row 8: BKT=BKT5, STATE=RAJASTHAN
$ws.Cells.Item(8, 1).Value = "BKT5"
$ws.Cells.Item(8, 2).Value = "RAJASTHAN"
$ws.Cells.Item(8, 3).Value = 2298733
$ws.Cells.Item(8, 4).Value = 2
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 2298733
$ws.Cells.Item(8, 14).Value = 0
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 0
$ws.Cells.Item(8, 18).Value = 0
$ws.Cells.Item(8, 19).Value = 0
$ws.Cells.Item(8, 20).Value = 0
$ws.Cells.Item(8, 21).Value = 100
$ws.Cells.Item(8, 22).Value = 0
$ws.Cells.Item(8, 23).Value = 0
$ws.Cells.Item(8, 24).Value = 0
$ws.Cells.Item(8, 25).Value = 0
$ws.Cells.Item(8, 26).Value = 0
$ws.Cells.Item(8, 27).Value = 0
$ws.Cells.Item(8, 28).Value = 0
$ws.Cells.Item(8, 29).Value = 0
$ws.Cells.Item(8, 30).Value = 0
$ws.Cells.Item(8, 31).Value = 0

# Row 9: BKT=BKT7, STATE=RAJASTHAN
$ws.Cells.Item(9, 1).Value = "BKT7"
$ws.Cells.Item(9, 2).Value = "RAJASTHAN"
$ws.Cells.Item(9, 3).Value = 1180938.14
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = 1180938.14
$ws.Cells.Item(9, 14).Value = 0
$ws.Cells.Item(9, 15).Value = 0
$ws.Cells.Item(9, 16).Value = 0
$ws.Cells.Item(9, 17).Value = 0
$ws.Cells.Item(9, 18).Value = 0
$ws.Cells.Item(9, 19).Value = 0
$ws.Cells.Item(9, 20).Value = 0
$ws.Cells.Item(9, 21).Value = 100
$ws.Cells.Item(9, 22).Value = 0
$ws.Cells.Item(9, 23).Value = 0
$ws.Cells.Item(9, 24).Value = 0
$ws.Cells.Item(9, 25).Value = 0
$ws.Cells.Item(9, 26).Value = 0
$ws.Cells.Item(9, 27).Value = 0
$ws.Cells.Item(9, 28).Value = 0
$ws.Cells.Item(9, 29).Value = 0
$ws.Cells.Item(9, 30).Value = 0
$ws.Cells.Item(9, 31).Value = 0

